$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.845.24'
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').Value = '1.629.80'
$ws.Range('E3').Value = '  -0.86%  '
$ws.Range('D4').Value = "'0.994"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.64%  '
$ws.Range('D5').Value = "'211.16"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.16%  '
$ws.Range('D6').Value = "'0.522"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.61%  '
$ws.Range('D7').Value = "'0.993"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.69%  '
$ws.Range('D8').Value = "'23.24"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.00%  '
$ws.Range('E9').Value = '  -2.63%  '
$ws.Range('D10').Value = "'0.0613"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.65%  '
$ws.Range('D11').Value = "'0.0878"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.06%  '
$ws.Range('D12').Value = '1.859.79'
$ws.Range('E12').Value = '  -0.92%  '
$ws.Range('D13').Value = '1.625.67'
$ws.Range('E13').Value = '  -1.11%  '
$ws.Range('D14').Value = "'4.04"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.77%  '
$ws.Range('D15').Value = "'0.562"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.88%  '
$ws.Range('D16').Value = "'65.33"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('D17').Value = '27.826.74'
$ws.Range('E17').Value = '  -0.27%  '
$ws.Range('D18').Value = "'231.84"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.51%  '
$ws.Range('D19').Value = '0.0₃0724'
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('D20').Value = "'7.53"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.37%  '
$ws.Range('D21').Value = "'0.994"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.68%  '
$ws.Range('D22').Value = "'10.40"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.86%  '
$ws.Range('D23').Value = "'4.34"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.21%  '
$ws.Range('D24').Value = "'2.05"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.81%  '
$ws.Range('D25').Value = "'153.81"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.72%  '
$ws.Range('D26').Value = "'6.89"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.42%  '
$ws.Range('D27').Value = "'0.111"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('D28').Value = "'15.61"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('D29').Value = "'0.994"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.64%  '
$ws.Range('D30').Value = "'1.18"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.21%  '
$ws.Range('D31').Value = "'0.0480"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.00%  '
$ws.Range('D32').Value = "'3.40"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.66%  '
$ws.Range('D33').Value = "'3.08"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.27%  '
$ws.Range('D34').Value = '1.399.93'
$ws.Range('E34').Value = '  -2.71%  '
$ws.Range('D35').Value = "'1.56"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.52%  '
$ws.Range('D36').Value = "'0.998"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.83%  '
$ws.Range('E37').Value = '  +0.41%  '
$ws.Range('D38').Value = "'0.0171"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.37%  '
$ws.Range('D39').Value = "'0.558"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('D40').Value = "'0.867"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.19%  '
$ws.Range('E41').Value = '  -1.43%  '
$ws.Range('D42').Value = "'0.993"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.70%  '
$ws.Range('D43').Value = "'66.64"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.88%  '
$ws.Range('D44').Value = "'1.83"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.33%  '
$ws.Range('D45').Value = "'5.44"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('E46').Value = '  -0.91%  '
$ws.Range('D47').Value = '1.769.01'
$ws.Range('E47').Value = '  -0.92%  '
$ws.Range('D48').Value = "'87.97"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.40%  '
$ws.Range('D49').Value = '0.0₆0102'
$ws.Range('E49').Value = '  -4.21%  '
$ws.Range('D50').Value = "'0.0998"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.05%  '
$ws.Range('D51').Value = "'0.0506"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.17%  '

"Update complete"
